# Updates cryptos list values (Price / Volume(1h) columns, plus the
# Cronos/EnergySwap row swap) to match the refreshed snapshot.
#
# The Price/Volume cells are plain text (e.g. "29.531.10", "  +0.31%  ")
# even though many look numeric. Excel auto-coerces a numeric-looking
# string assigned via .Value into a real number, so we lead each new
# value with a literal apostrophe to force text, then immediately reset
# the cell style back to "Normal" so no stray Text number-format / style
# is left behind (matches the original, unstyled cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2, E2
$ws.Range("D2").Value = "'29.501.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.23%  "
$ws.Range("E2").Style = "Normal"

# Row 3: D3, E3
$ws.Range("D3").Value = "'1.904.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.62%  "
$ws.Range("E3").Style = "Normal"

# Row 4: E4
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"

# Row 5: D5, E5
$ws.Range("D5").Value = "'337.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.05%  "
$ws.Range("E5").Style = "Normal"

# Row 6: D6, E6
$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.02%  "
$ws.Range("E6").Style = "Normal"

# Row 7: D7, E7
$ws.Range("D7").Value = "'0.4762"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.15%  "
$ws.Range("E7").Style = "Normal"

# Row 8: D8, E8
$ws.Range("D8").Value = "'0.4001"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.76%  "
$ws.Range("E8").Style = "Normal"

# Row 9: D9, E9
$ws.Range("D9").Value = "'0.08031"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.45%  "
$ws.Range("E9").Style = "Normal"

# Row 10: D10
$ws.Range("D10").Value = "'0.9911"
$ws.Range("D10").Style = "Normal"

# Row 11: D11, E11
$ws.Range("D11").Value = "'23.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.76%  "
$ws.Range("E11").Style = "Normal"

# Row 12: D12, E12
$ws.Range("D12").Value = "'1.902.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.27%  "
$ws.Range("E12").Style = "Normal"

# Row 13: D13, E13
$ws.Range("D13").Value = "'5.913"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.61%  "
$ws.Range("E13").Style = "Normal"

# Row 14: D14, E14
$ws.Range("D14").Value = "'7.106"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.86%  "
$ws.Range("E14").Style = "Normal"

# Row 15: E15
$ws.Range("E15").Value = "'  -2.53%  "
$ws.Range("E15").Style = "Normal"

# Row 16: D16, E16
$ws.Range("D16").Value = "'0.06830"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.30%  "
$ws.Range("E16").Style = "Normal"

# Row 17: D17
$ws.Range("D17").Value = "'1.008"
$ws.Range("D17").Style = "Normal"

# Row 18: D18, E18
$ws.Range("D18").Value = "'0.00001020"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.88%  "
$ws.Range("E18").Style = "Normal"

# Row 19: D19, E19
$ws.Range("D19").Value = "'17.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.94%  "
$ws.Range("E19").Style = "Normal"

# Row 20: D20, E20
$ws.Range("D20").Value = "'1.006"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.13%  "
$ws.Range("E20").Style = "Normal"

# Row 21: D21, E21
$ws.Range("D21").Value = "'29.513.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.22%  "
$ws.Range("E21").Style = "Normal"

# Row 22: D22, E22
$ws.Range("D22").Value = "'5.501"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.73%  "
$ws.Range("E22").Style = "Normal"

# Row 23: D23, E23
$ws.Range("D23").Value = "'11.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.93%  "
$ws.Range("E23").Style = "Normal"

# Row 24: D24, E24
$ws.Range("D24").Value = "'2.156"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.96%  "
$ws.Range("E24").Style = "Normal"

# Row 25: D25, E25
$ws.Range("D25").Value = "'2.138.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.35%  "
$ws.Range("E25").Style = "Normal"

# Row 26: D26, E26
$ws.Range("D26").Value = "'156.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.92%  "
$ws.Range("E26").Style = "Normal"

# Row 27: D27, E27
$ws.Range("D27").Value = "'6.464"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.10%  "
$ws.Range("E27").Style = "Normal"

# Row 28: D28, E28
$ws.Range("D28").Value = "'19.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.67%  "
$ws.Range("E28").Style = "Normal"

# Row 29: D29, E29
$ws.Range("D29").Value = "'2.050"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.97%  "
$ws.Range("E29").Style = "Normal"

# Row 30: D30, E30
$ws.Range("D30").Value = "'119.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.23%  "
$ws.Range("E30").Style = "Normal"

# Row 31: D31, E31
$ws.Range("D31").Value = "'0.9946"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.43%  "
$ws.Range("E31").Style = "Normal"

# Row 32: D32, E32
$ws.Range("D32").Value = "'0.09523"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.73%  "
$ws.Range("E32").Style = "Normal"

# Row 33: D33, E33
$ws.Range("D33").Value = "'5.466"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.90%  "
$ws.Range("E33").Style = "Normal"

# Row 34: D34, E34
$ws.Range("D34").Value = "'3.540"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.40%  "
$ws.Range("E34").Style = "Normal"

# Row 35: D35, E35
$ws.Range("D35").Value = "'1.385"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.98%  "
$ws.Range("E35").Style = "Normal"

# Row 36: D36, E36
$ws.Range("D36").Value = "'0.06455"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.63%  "
$ws.Range("E36").Style = "Normal"

# Row 37: D37, E37
$ws.Range("D37").Value = "'0.02238"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.00%  "
$ws.Range("E37").Style = "Normal"

# Row 38: D38, E38
$ws.Range("D38").Value = "'1.190"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.88%  "
$ws.Range("E38").Style = "Normal"

# Row 39: D39, E39
$ws.Range("D39").Value = "'0.5815"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.77%  "
$ws.Range("E39").Style = "Normal"

# Row 40: D40, E40
$ws.Range("D40").Value = "'10.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.80%  "
$ws.Range("E40").Style = "Normal"

# Row 41: D41, E41
$ws.Range("D41").Value = "'7.751"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.82%  "
$ws.Range("E41").Style = "Normal"

# Row 42: D42, E42
$ws.Range("D42").Value = "'0.1820"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.58%  "
$ws.Range("E42").Style = "Normal"

# Row 43: D43, E43
$ws.Range("D43").Value = "'2.452"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.06%  "
$ws.Range("E43").Style = "Normal"

# Row 44: D44, E44
$ws.Range("D44").Value = "'1.235"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.45%  "
$ws.Range("E44").Style = "Normal"

# Row 45: B45, C45, D45, E45
$ws.Range("B45").Value = "'EnergySwap"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'12.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.44%  "
$ws.Range("E45").Style = "Normal"

# Row 46: B46, C46, D46, E46
$ws.Range("B46").Value = "'Cronos"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.07416"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.45%  "
$ws.Range("E46").Style = "Normal"

# Row 47: D47, E47
$ws.Range("D47").Value = "'0.5478"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.02%  "
$ws.Range("E47").Style = "Normal"

# Row 48: D48, E48
$ws.Range("D48").Value = "'1.945"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.51%  "
$ws.Range("E48").Style = "Normal"

# Row 49: D49, E49
$ws.Range("D49").Value = "'115.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.95%  "
$ws.Range("E49").Style = "Normal"

# Row 50: D50, E50
$ws.Range("D50").Value = "'2.378"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.17%  "
$ws.Range("E50").Style = "Normal"

# Row 51: D51, E51
$ws.Range("D51").Value = "'71.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.77%  "
$ws.Range("E51").Style = "Normal"
